$wb = $excel.ActiveWorkbook

# ---- Sheet1: "Schedule" ----
$wsSchedule = $wb.Worksheets.Item("Schedule")

$arr1 = New-Object 'object[,]' 4,6
$arr1[0,0] = 46037
$arr1[0,1] = 46037.20833333334
$arr1[0,2] = 5
$arr1[0,3] = 18.9
$arr1[0,4] = 737.8842217499999
$arr1[0,5] = 39.04149321428572
$arr1[1,0] = 46037.29166666666
$arr1[1,1] = 46037.66666666666
$arr1[1,2] = 9
$arr1[1,3] = 34.02
$arr1[1,4] = 583.4340525
$arr1[1,5] = 17.14973699294533
$arr1[2,0] = 46037.91666666666
$arr1[2,1] = 46038.125
$arr1[2,2] = 5
$arr1[2,3] = 18.9
$arr1[2,4] = 532.4608867500001
$arr1[2,5] = 28.1725336904762
$arr1[3,0] = 46038.29166666666
$arr1[3,1] = 46038.66666666666
$arr1[3,2] = 9
$arr1[3,3] = 34.02
$arr1[3,4] = 483.8738970000001
$arr1[3,5] = 14.22321860670194

$wsSchedule.Range("A2:F5").Value = $arr1

# Start/Stop Time columns use the workbook's existing DateTime display
# format (same style already applied to A2/B2) - re-apply across the
# whole column so the newly-added rows 3-5 match it too.
$wsSchedule.Range("A2:B5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---- Sheet2: "Detailed" ----
$wsDetailed = $wb.Worksheets.Item("Detailed")

$arr2 = New-Object 'object[,]' 96,5
$arr2[0,0] = 46037
$arr2[0,1] = 78
$arr2[0,2] = "historical"
$arr2[0,3] = 46037
$arr2[0,4] = "ON"
$arr2[1,0] = 46037.02083333334
$arr2[1,1] = 76.66624
$arr2[1,2] = "historical"
$arr2[1,3] = 46037
$arr2[1,4] = "ON"
$arr2[2,0] = 46037.04166666666
$arr2[2,1] = 78
$arr2[2,2] = "historical"
$arr2[2,3] = 46037
$arr2[2,4] = "ON"
$arr2[3,0] = 46037.0625
$arr2[3,1] = 78
$arr2[3,2] = "historical"
$arr2[3,3] = 46037
$arr2[3,4] = "ON"
$arr2[4,0] = 46037.08333333334
$arr2[4,1] = 78
$arr2[4,2] = "historical"
$arr2[4,3] = 46037
$arr2[4,4] = "ON"
$arr2[5,0] = 46037.10416666666
$arr2[5,1] = 78
$arr2[5,2] = "historical"
$arr2[5,3] = 46037
$arr2[5,4] = "ON"
$arr2[6,0] = 46037.125
$arr2[6,1] = 78
$arr2[6,2] = "historical"
$arr2[6,3] = 46037
$arr2[6,4] = "ON"
$arr2[7,0] = 46037.14583333334
$arr2[7,1] = 73.78127000000001
$arr2[7,2] = "historical"
$arr2[7,3] = 46037
$arr2[7,4] = "ON"
$arr2[8,0] = 46037.16666666666
$arr2[8,1] = 73.1985
$arr2[8,2] = "historical"
$arr2[8,3] = 46037
$arr2[8,4] = "ON"
$arr2[9,0] = 46037.1875
$arr2[9,1] = 65.15832
$arr2[9,2] = "historical"
$arr2[9,3] = 46037
$arr2[9,4] = "ON"
$arr2[10,0] = 46037.20833333334
$arr2[10,1] = 81.07834
$arr2[10,2] = "historical"
$arr2[10,3] = 46037
$arr2[10,4] = "OFF"
$arr2[11,0] = 46037.22916666666
$arr2[11,1] = 93.42456
$arr2[11,2] = "historical"
$arr2[11,3] = 46037
$arr2[11,4] = "OFF"
$arr2[12,0] = 46037.25
$arr2[12,1] = 90.79733
$arr2[12,2] = "historical"
$arr2[12,3] = 46037
$arr2[12,4] = "OFF"
$arr2[13,0] = 46037.27083333334
$arr2[13,1] = 79.09635
$arr2[13,2] = "historical"
$arr2[13,3] = 46037
$arr2[13,4] = "OFF"
$arr2[14,0] = 46037.29166666666
$arr2[14,1] = 56.98
$arr2[14,2] = "historical"
$arr2[14,3] = 46037
$arr2[14,4] = "ON"
$arr2[15,0] = 46037.3125
$arr2[15,1] = 51.40072
$arr2[15,2] = "historical"
$arr2[15,3] = 46037
$arr2[15,4] = "ON"
$arr2[16,0] = 46037.33333333334
$arr2[16,1] = 50.14948
$arr2[16,2] = "historical"
$arr2[16,3] = 46037
$arr2[16,4] = "ON"
$arr2[17,0] = 46037.35416666666
$arr2[17,1] = 36.06
$arr2[17,2] = "historical"
$arr2[17,3] = 46037
$arr2[17,4] = "ON"
$arr2[18,0] = 46037.375
$arr2[18,1] = 36.06
$arr2[18,2] = "historical"
$arr2[18,3] = 46037
$arr2[18,4] = "ON"
$arr2[19,0] = 46037.39583333334
$arr2[19,1] = 8.18614
$arr2[19,2] = "historical"
$arr2[19,3] = 46037
$arr2[19,4] = "ON"
$arr2[20,0] = 46037.41666666666
$arr2[20,1] = 10.3824
$arr2[20,2] = "historical"
$arr2[20,3] = 46037
$arr2[20,4] = "ON"
$arr2[21,0] = 46037.4375
$arr2[21,1] = 36.06
$arr2[21,2] = "historical"
$arr2[21,3] = 46037
$arr2[21,4] = "ON"
$arr2[22,0] = 46037.45833333334
$arr2[22,1] = 0.51
$arr2[22,2] = "historical"
$arr2[22,3] = 46037
$arr2[22,4] = "ON"
$arr2[23,0] = 46037.47916666666
$arr2[23,1] = 0.51
$arr2[23,2] = "historical"
$arr2[23,3] = 46037
$arr2[23,4] = "ON"
$arr2[24,0] = 46037.5
$arr2[24,1] = 34.01
$arr2[24,2] = "historical"
$arr2[24,3] = 46037
$arr2[24,4] = "ON"
$arr2[25,0] = 46037.52083333334
$arr2[25,1] = 36.0601
$arr2[25,2] = "historical"
$arr2[25,3] = 46037
$arr2[25,4] = "ON"
$arr2[26,0] = 46037.54166666666
$arr2[26,1] = 0.51
$arr2[26,2] = "historical"
$arr2[26,3] = 46037
$arr2[26,4] = "ON"
$arr2[27,0] = 46037.5625
$arr2[27,1] = 39.0601
$arr2[27,2] = "historical"
$arr2[27,3] = 46037
$arr2[27,4] = "ON"
$arr2[28,0] = 46037.58333333334
$arr2[28,1] = 52.32496
$arr2[28,2] = "historical"
$arr2[28,3] = 46037
$arr2[28,4] = "ON"
$arr2[29,0] = 46037.60416666666
$arr2[29,1] = 57.08
$arr2[29,2] = "forecast"
$arr2[29,3] = 46037
$arr2[29,4] = "ON"
$arr2[30,0] = 46037.625
$arr2[30,1] = 36.07
$arr2[30,2] = "forecast"
$arr2[30,3] = 46037
$arr2[30,4] = "ON"
$arr2[31,0] = 46037.64583333334
$arr2[31,1] = 56.98
$arr2[31,2] = "forecast"
$arr2[31,3] = 46037
$arr2[31,4] = "ON"
$arr2[32,0] = 46037.66666666666
$arr2[32,1] = 44.01769
$arr2[32,2] = "forecast"
$arr2[32,3] = 46037
$arr2[32,4] = "OFF"
$arr2[33,0] = 46037.6875
$arr2[33,1] = 49.15376
$arr2[33,2] = "forecast"
$arr2[33,3] = 46037
$arr2[33,4] = "OFF"
$arr2[34,0] = 46037.70833333334
$arr2[34,1] = 54.3948
$arr2[34,2] = "forecast"
$arr2[34,3] = 46037
$arr2[34,4] = "OFF"
$arr2[35,0] = 46037.72916666666
$arr2[35,1] = 18.54764
$arr2[35,2] = "forecast"
$arr2[35,3] = 46037
$arr2[35,4] = "OFF"
$arr2[36,0] = 46037.75
$arr2[36,1] = 55.11462
$arr2[36,2] = "forecast"
$arr2[36,3] = 46037
$arr2[36,4] = "OFF"
$arr2[37,0] = 46037.77083333334
$arr2[37,1] = 70.47145
$arr2[37,2] = "forecast"
$arr2[37,3] = 46037
$arr2[37,4] = "OFF"
$arr2[38,0] = 46037.79166666666
$arr2[38,1] = 120.01
$arr2[38,2] = "forecast"
$arr2[38,3] = 46037
$arr2[38,4] = "OFF"
$arr2[39,0] = 46037.8125
$arr2[39,1] = 120.01
$arr2[39,2] = "forecast"
$arr2[39,3] = 46037
$arr2[39,4] = "OFF"
$arr2[40,0] = 46037.83333333334
$arr2[40,1] = 120.01
$arr2[40,2] = "forecast"
$arr2[40,3] = 46037
$arr2[40,4] = "OFF"
$arr2[41,0] = 46037.85416666666
$arr2[41,1] = 101.25
$arr2[41,2] = "forecast"
$arr2[41,3] = 46037
$arr2[41,4] = "OFF"
$arr2[42,0] = 46037.875
$arr2[42,1] = 85.95
$arr2[42,2] = "forecast"
$arr2[42,3] = 46037
$arr2[42,4] = "OFF"
$arr2[43,0] = 46037.89583333334
$arr2[43,1] = 83.95355000000001
$arr2[43,2] = "forecast"
$arr2[43,3] = 46037
$arr2[43,4] = "OFF"
$arr2[44,0] = 46037.91666666666
$arr2[44,1] = 69.09249
$arr2[44,2] = "forecast"
$arr2[44,3] = 46037
$arr2[44,4] = "ON"
$arr2[45,0] = 46037.9375
$arr2[45,1] = 57.09
$arr2[45,2] = "forecast"
$arr2[45,3] = 46037
$arr2[45,4] = "ON"
$arr2[46,0] = 46037.95833333334
$arr2[46,1] = 57.09
$arr2[46,2] = "forecast"
$arr2[46,3] = 46037
$arr2[46,4] = "ON"
$arr2[47,0] = 46037.97916666666
$arr2[47,1] = 57.41519
$arr2[47,2] = "forecast"
$arr2[47,3] = 46037
$arr2[47,4] = "ON"
$arr2[48,0] = 46038
$arr2[48,1] = 57.06008
$arr2[48,2] = "forecast"
$arr2[48,3] = 46038
$arr2[48,4] = "ON"
$arr2[49,0] = 46038.02083333334
$arr2[49,1] = 57.06003
$arr2[49,2] = "forecast"
$arr2[49,3] = 46038
$arr2[49,4] = "ON"
$arr2[50,0] = 46038.04166666666
$arr2[50,1] = 56.98
$arr2[50,2] = "forecast"
$arr2[50,3] = 46038
$arr2[50,4] = "ON"
$arr2[51,0] = 46038.0625
$arr2[51,1] = 49.66065
$arr2[51,2] = "forecast"
$arr2[51,3] = 46038
$arr2[51,4] = "ON"
$arr2[52,0] = 46038.08333333334
$arr2[52,1] = 48.59529
$arr2[52,2] = "forecast"
$arr2[52,3] = 46038
$arr2[52,4] = "ON"
$arr2[53,0] = 46038.10416666666
$arr2[53,1] = 36.07
$arr2[53,2] = "forecast"
$arr2[53,3] = 46038
$arr2[53,4] = "ON"
$arr2[54,0] = 46038.125
$arr2[54,1] = 36.07
$arr2[54,2] = "forecast"
$arr2[54,3] = 46038
$arr2[54,4] = "OFF"
$arr2[55,0] = 46038.14583333334
$arr2[55,1] = 42.26032
$arr2[55,2] = "forecast"
$arr2[55,3] = 46038
$arr2[55,4] = "OFF"
$arr2[56,0] = 46038.16666666666
$arr2[56,1] = 55.39499
$arr2[56,2] = "forecast"
$arr2[56,3] = 46038
$arr2[56,4] = "OFF"
$arr2[57,0] = 46038.1875
$arr2[57,1] = 57.05991
$arr2[57,2] = "forecast"
$arr2[57,3] = 46038
$arr2[57,4] = "OFF"
$arr2[58,0] = 46038.20833333334
$arr2[58,1] = 57.06003
$arr2[58,2] = "forecast"
$arr2[58,3] = 46038
$arr2[58,4] = "OFF"
$arr2[59,0] = 46038.22916666666
$arr2[59,1] = 60.57749
$arr2[59,2] = "forecast"
$arr2[59,3] = 46038
$arr2[59,4] = "OFF"
$arr2[60,0] = 46038.25
$arr2[60,1] = 57.07828
$arr2[60,2] = "forecast"
$arr2[60,3] = 46038
$arr2[60,4] = "OFF"
$arr2[61,0] = 46038.27083333334
$arr2[61,1] = 56.98
$arr2[61,2] = "forecast"
$arr2[61,3] = 46038
$arr2[61,4] = "OFF"
$arr2[62,0] = 46038.29166666666
$arr2[62,1] = 36.05972
$arr2[62,2] = "forecast"
$arr2[62,3] = 46038
$arr2[62,4] = "ON"
$arr2[63,0] = 46038.3125
$arr2[63,1] = 36.06
$arr2[63,2] = "forecast"
$arr2[63,3] = 46038
$arr2[63,4] = "ON"
$arr2[64,0] = 46038.33333333334
$arr2[64,1] = 36.06
$arr2[64,2] = "forecast"
$arr2[64,3] = 46038
$arr2[64,4] = "ON"
$arr2[65,0] = 46038.35416666666
$arr2[65,1] = 36.06
$arr2[65,2] = "forecast"
$arr2[65,3] = 46038
$arr2[65,4] = "ON"
$arr2[66,0] = 46038.375
$arr2[66,1] = 36.06
$arr2[66,2] = "forecast"
$arr2[66,3] = 46038
$arr2[66,4] = "ON"
$arr2[67,0] = 46038.39583333334
$arr2[67,1] = 22.62945
$arr2[67,2] = "forecast"
$arr2[67,3] = 46038
$arr2[67,4] = "ON"
$arr2[68,0] = 46038.41666666666
$arr2[68,1] = 36.06
$arr2[68,2] = "forecast"
$arr2[68,3] = 46038
$arr2[68,4] = "ON"
$arr2[69,0] = 46038.4375
$arr2[69,1] = 36.07
$arr2[69,2] = "forecast"
$arr2[69,3] = 46038
$arr2[69,4] = "ON"
$arr2[70,0] = 46038.45833333334
$arr2[70,1] = 36.05949
$arr2[70,2] = "forecast"
$arr2[70,3] = 46038
$arr2[70,4] = "ON"
$arr2[71,0] = 46038.47916666666
$arr2[71,1] = 36.05989
$arr2[71,2] = "forecast"
$arr2[71,3] = 46038
$arr2[71,4] = "ON"
$arr2[72,0] = 46038.5
$arr2[72,1] = 36.06
$arr2[72,2] = "forecast"
$arr2[72,3] = 46038
$arr2[72,4] = "ON"
$arr2[73,0] = 46038.52083333334
$arr2[73,1] = 36.0601
$arr2[73,2] = "forecast"
$arr2[73,3] = 46038
$arr2[73,4] = "ON"
$arr2[74,0] = 46038.54166666666
$arr2[74,1] = 36.0601
$arr2[74,2] = "forecast"
$arr2[74,3] = 46038
$arr2[74,4] = "ON"
$arr2[75,0] = 46038.5625
$arr2[75,1] = 22.07
$arr2[75,2] = "forecast"
$arr2[75,3] = 46038
$arr2[75,4] = "ON"
$arr2[76,0] = 46038.58333333334
$arr2[76,1] = -7.01
$arr2[76,2] = "forecast"
$arr2[76,3] = 46038
$arr2[76,4] = "ON"
$arr2[77,0] = 46038.60416666666
$arr2[77,1] = -2.47963
$arr2[77,2] = "forecast"
$arr2[77,3] = 46038
$arr2[77,4] = "ON"
$arr2[78,0] = 46038.625
$arr2[78,1] = 0.57069
$arr2[78,2] = "forecast"
$arr2[78,3] = 46038
$arr2[78,4] = "ON"
$arr2[79,0] = 46038.64583333334
$arr2[79,1] = 27.77111
$arr2[79,2] = "forecast"
$arr2[79,3] = 46038
$arr2[79,4] = "ON"
$arr2[80,0] = 46038.66666666666
$arr2[80,1] = 27.7711
$arr2[80,2] = "forecast"
$arr2[80,3] = 46038
$arr2[80,4] = "OFF"
$arr2[81,0] = 46038.6875
$arr2[81,1] = 21.4936
$arr2[81,2] = "forecast"
$arr2[81,3] = 46038
$arr2[81,4] = "OFF"
$arr2[82,0] = 46038.70833333334
$arr2[82,1] = 43.24919
$arr2[82,2] = "forecast"
$arr2[82,3] = 46038
$arr2[82,4] = "OFF"
$arr2[83,0] = 46038.72916666666
$arr2[83,1] = 37.83416
$arr2[83,2] = "forecast"
$arr2[83,3] = 46038
$arr2[83,4] = "OFF"
$arr2[84,0] = 46038.75
$arr2[84,1] = 0.48373
$arr2[84,2] = "forecast"
$arr2[84,3] = 46038
$arr2[84,4] = "OFF"
$arr2[85,0] = 46038.77083333334
$arr2[85,1] = 53.90481
$arr2[85,2] = "forecast"
$arr2[85,3] = 46038
$arr2[85,4] = "OFF"
$arr2[86,0] = 46038.79166666666
$arr2[86,1] = 299.98
$arr2[86,2] = "forecast"
$arr2[86,3] = 46038
$arr2[86,4] = "OFF"
$arr2[87,0] = 46038.8125
$arr2[87,1] = 299.98
$arr2[87,2] = "forecast"
$arr2[87,3] = 46038
$arr2[87,4] = "OFF"
$arr2[88,0] = 46038.83333333334
$arr2[88,1] = 67.39879999999999
$arr2[88,2] = "forecast"
$arr2[88,3] = 46038
$arr2[88,4] = "OFF"
$arr2[89,0] = 46038.85416666666
$arr2[89,1] = 49.23153
$arr2[89,2] = "forecast"
$arr2[89,3] = 46038
$arr2[89,4] = "OFF"
$arr2[90,0] = 46038.875
$arr2[90,1] = 40.98924
$arr2[90,2] = "forecast"
$arr2[90,3] = 46038
$arr2[90,4] = "OFF"
$arr2[91,0] = 46038.89583333334
$arr2[91,1] = 40.5543
$arr2[91,2] = "forecast"
$arr2[91,3] = 46038
$arr2[91,4] = "OFF"
$arr2[92,0] = 46038.91666666666
$arr2[92,1] = 40.56485
$arr2[92,2] = "forecast"
$arr2[92,3] = 46038
$arr2[92,4] = "OFF"
$arr2[93,0] = 46038.9375
$arr2[93,1] = 48.38127
$arr2[93,2] = "forecast"
$arr2[93,3] = 46038
$arr2[93,4] = "OFF"
$arr2[94,0] = 46038.95833333334
$arr2[94,1] = 47.88557
$arr2[94,2] = "forecast"
$arr2[94,3] = 46038
$arr2[94,4] = "OFF"
$arr2[95,0] = 46038.97916666666
$arr2[95,1] = 40.46757
$arr2[95,2] = "forecast"
$arr2[95,3] = 46038
$arr2[95,4] = "OFF"

$wsDetailed.Range("A2:E97").Value = $arr2

# Re-apply the DateTime / Date number formats so the newly-created rows
# (50-97) pick up the same display format already used by rows 2-49,
# matching the existing style definitions (numFmtId 165 / 167) rather than
# creating new ones.
$wsDetailed.Range("A2:A97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsDetailed.Range("D2:D97").NumberFormat = "YYYY-MM-DD"
